# Work on line graph
# Adds a new "Line Graph:" labelled block (header row + data rows) below the
# existing Table1, extends the table by one row to cover the new label row,
# updates the sheet dimension/selection, and sets the page to portrait.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Grow Table1 so it also covers row 38 (the new "Line Graph:" label row).
# ---------------------------------------------------------------------------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:L38"))

# ---------------------------------------------------------------------------
# 2. Bold "Line Graph:" label in A38.
# ---------------------------------------------------------------------------
$ws.Range("A38").Value = "Line Graph:"
$ws.Range("A38").Font.Bold = $true

# ---------------------------------------------------------------------------
# 3. New header + data block in rows 39-63 (outside of Table1).
# ---------------------------------------------------------------------------
$data = @(
    @("dataset", "embedding_method", "reference_vector", "euclidian", "cosine", "minkowski", "mahalanobis", "braycurtis", "canberra", "chebyshev", "cityblock", "correlation"),
    @("RMV", "Graph2Vec", "mean", 0.44, 0.44, 0.44, 0.25, 0.56, 0.63, 0.44, 0.44, 0.44),
    @("RMV", "Graph2Vec", "initial", 0.45, 0.35, 0.46, 0.0, 0.3, 0.32, 0.43, 0.46, 0.35),
    @("RMV", "LDP", "mean", 0.27, 0.21, 0.32, 0.0, 0.23, 0.23, 0.21, 0.32, 0.21),
    @("RMV", "LDP", "initial", 0.2, 0.24, 0.21, 0.0, 0.2, 0.22, 0.2, 0.21, 0.24),
    @("RMV", "FeatherGraph", "mean", 0.22, 0.21, 0.23, 0.0, 0.23, 0.22, 0.22, 0.23, 0.22),
    @("RMV", "FeatherGraph", "initial", 0.24, 0.25, 0.25, 0.0, 0.27, 0.23, 0.24, 0.25, 0.24),
    @("RMS", "Graph2Vec", "mean", 0.29, 0.39, 0.28, 0.28, 0.26, 0.28, 0.3, 0.28, 0.39),
    @("RMS", "Graph2Vec", "initial", 0.26, 0.46, 0.27, 0.0, 0.23, 0.22, 0.22, 0.27, 0.46),
    @("RMS", "LDP", "mean", 0.36, 0.2, 0.37, 0.0, 0.22, 0.22, 0.36, 0.37, 0.21),
    @("RMS", "LDP", "initial", 0.2, 0.21, 0.2, 0.0, 0.2, 0.21, 0.2, 0.2, 0.21),
    @("RMS", "FeatherGraph", "mean", 0.25, 0.26, 0.25, 0.0, 0.25, 0.24, 0.25, 0.25, 0.26),
    @("RMS", "FeatherGraph", "initial", 0.25, 0.26, 0.26, 0.0, 0.27, 0.25, 0.25, 0.26, 0.26),
    @("ENRON", "Graph2Vec", "mean", 0.07, 0.1, 0.07, 0.05, 0.07, 0.08, 0.07, 0.07, 0.1),
    @("ENRON", "Graph2Vec", "initial", 0.07, 0.1, 0.07, 0.0, 0.07, 0.07, 0.07, 0.07, 0.1),
    @("ENRON", "LDP", "mean", 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0),
    @("ENRON", "LDP", "initial", 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0),
    @("ENRON", "FeatherGraph", "mean", 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0),
    @("ENRON", "FeatherGraph", "initial", 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0, 0.0),
    @("CS", "Graph2Vec", "mean", 0.42, 0.42, 0.38, 0.67, 0.48, 0.4, 0.63, 0.38, 0.4),
    @("CS", "Graph2Vec", "initial", 0.44, 0.45, 0.44, 0.0, 0.48, 0.5, 0.44, 0.44, 0.49),
    @("CS", "LDP", "mean", 0.48, 0.45, 0.5, 0.0, 0.48, 0.55, 0.48, 0.5, 0.6),
    @("CS", "LDP", "initial", 0.59, 0.6, 0.59, 0.0, 0.59, 0.52, 0.59, 0.59, 0.51),
    @("CS", "FeatherGraph", "mean", 0.54, 0.63, 0.54, 0.0, 0.55, 0.49, 0.49, 0.54, 0.5),
    @("CS", "FeatherGraph", "initial", 0.54, 0.5, 0.54, 0.4, 0.54, 0.57, 0.5, 0.54, 0.5)
)

$startRow = 39
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    $arr = New-Object 'object[,]' 1,12
    for ($j = 0; $j -lt 12; $j++) {
        $arr[0,$j] = $rowVals[$j]
    }
    $r = $startRow + $i
    $rng = $ws.Range("A" + $r + ":L" + $r)
    $rng.Value2 = $arr
}

# ---------------------------------------------------------------------------
# 4. Page setup: portrait orientation.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 5. Update selection / view to reflect where the user ended up editing.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("G69").Select()
